$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7175764140682759
$ws.Range("C2").Value = 0.06260135146008849
$ws.Range("D2").Value = 0.1265557049407491
$ws.Range("E2").Value = 0.1215562336815807
$ws.Range("F2").Value = 1.730212186002525
$ws.Range("I2").Value = 1.091659938521605
$ws.Range("J2").Value = 0.1506897988592826
$ws.Range("K2").Value = 0.442376835111105
$ws.Range("L2").Value = 0.3010045715382859
$ws.Range("N2").Value = 2.031795384352552
$ws.Range("O2").Value = 4.502243148714626
$ws.Range("B3").Value = 0.680771382272269
$ws.Range("C3").Value = 0.0603773326943724
$ws.Range("D3").Value = 0.1246814889621248
$ws.Range("E3").Value = 0.1215706998981005
$ws.Range("F3").Value = 1.736667472839201
$ws.Range("I3").Value = 1.099026108219643
$ws.Range("J3").Value = 0.1515007563693196
$ws.Range("K3").Value = 0.4075864759762453
$ws.Range("L3").Value = 0.296888353869754
$ws.Range("N3").Value = 2.050031960776172
$ws.Range("O3").Value = 4.526135960436591
$ws.Range("B4").Value = 0.6583985860095254
$ws.Range("C4").Value = 0.05899572011686871
$ws.Range("D4").Value = 0.1235750177935913
$ws.Range("E4").Value = 0.1216189954940319
$ws.Range("F4").Value = 1.741356316376994
$ws.Range("I4").Value = 1.103977940612751
$ws.Range("J4").Value = 0.1520419333065863
$ws.Range("K4").Value = 0.386321024341413
$ws.Range("L4").Value = 0.2944776136466842
$ws.Range("N4").Value = 2.061805402870696
$ws.Range("O4").Value = 4.54271532295482
$ws.Range("B5").Value = 0.6493390029859256
$ws.Range("C5").Value = 0.05842868181383665
$ws.Range("D5").Value = 0.1231353322108646
$ws.Range("E5").Value = 0.1216486175234266
$ws.Range("F5").Value = 1.743449680609771
$ws.Range("I5").Value = 1.10610382385719
$ws.Range("J5").Value = 0.1522733612966842
$ws.Range("K5").Value = 0.3776799283747891
$ws.Range("L5").Value = 0.2935246841258845
$ws.Range("N5").Value = 2.066748110782149
$ws.Range("O5").Value = 4.549952074521983
$ws.Range("B6").Value = 0.6478381625347254
$ws.Range("C6").Value = 0.05833428325504286
$ws.Range("D6").Value = 0.1230630018806949
$ws.Range("E6").Value = 0.1216541375553106
$ws.Range("F6").Value = 1.743808319357782
$ws.Range("I6").Value = 1.106463349363707
$ws.Range("J6").Value = 0.1523124482826823
$ws.Range("K6").Value = 0.3762465938263944
$ws.Range("L6").Value = 0.2933682346817932
$ws.Range("N6").Value = 2.06757759936855
$ws.Range("O6").Value = 4.551182768170293
$ws.Range("B7").Value = 0.658276171410705
$ws.Range("C7").Value = 0.05898808909299191
$ws.Range("D7").Value = 0.1235690425567029
$ws.Range("E7").Value = 0.1216193546950173
$ws.Range("F7").Value = 1.741383808459915
$ws.Range("I7").Value = 1.104006173710069
$ws.Range("J7").Value = 0.1520450102900242
$ws.Range("K7").Value = 0.3862043864783686
$ws.Range("L7").Value = 0.2944646426122617
$ws.Range("N7").Value = 2.061871475033595
$ws.Range("O7").Value = 4.542810974055385
$ws.Range("B8").Value = 0.7048396587177592
$ws.Range("C8").Value = 0.06183785212275694
$ws.Range("D8").Value = 0.125900320674603
$ws.Range("E8").Value = 0.1215530594235261
$ws.Range("F8").Value = 1.732287585112459
$ws.Range("I8").Value = 1.094110821374876
$ws.Range("J8").Value = 0.1509604530448243
$ws.Range("K8").Value = 0.4303615592287429
$ws.Range("L8").Value = 0.2995611732831662
$ws.Range("N8").Value = 2.03796382606891
$ws.Range("O8").Value = 4.510085441440353
$ws.Range("B9").Value = 0.7979118093734314
$ws.Range("C9").Value = 0.06729829495397155
$ws.Range("D9").Value = 0.1308207564177053
$ws.Range("E9").Value = 0.1217345333303683
$ws.Range("F9").Value = 1.720194683037604
$ws.Range("I9").Value = 1.078105420177238
$ws.Range("J9").Value = 0.1491759215874673
$ws.Range("K9").Value = 0.5176924347111651
$ws.Range("L9").Value = 0.3104755802110475
$ws.Range("N9").Value = 1.995648494263387
$ws.Range("O9").Value = 4.461039308230681
$ws.Range("B10").Value = 0.8673338054271085
$ws.Range("C10").Value = 0.07123171116337801
$ws.Range("D10").Value = 0.1346452404174556
$ws.Range("E10").Value = 0.122056198820232
$ws.Range("F10").Value = 1.714800182459754
$ws.Range("I10").Value = 1.068412909307192
$ws.Range("J10").Value = 0.1480723474264281
$ws.Range("K10").Value = 0.5822807839025472
$ws.Range("L10").Value = 0.3190492983363384
$ws.Range("N10").Value = 1.96733784435899
$ws.Range("O10").Value = 4.434204547313129
$ws.Range("B11").Value = 0.8991355763765512
$ws.Range("C11").Value = 0.07300404535126859
$ws.Range("D11").Value = 0.136429924199831
$ws.Range("E11").Value = 0.1222431084518938
$ws.Range("F11").Value = 1.713101494637058
$ws.Range("I11").Value = 1.064451115026799
$ws.Range("J11").Value = 0.1476151241628365
$ws.Range("K11").Value = 0.611751552559781
$ws.Range("L11").Value = 0.3230689841359151
$ws.Range("N11").Value = 1.955060619126506
$ws.Range("O11").Value = 4.423989558507401
$ws.Range("B12").Value = 0.9112092156240124
$ws.Range("C12").Value = 0.07367272454864349
$ws.Range("D12").Value = 0.1371121303970568
$ws.Range("E12").Value = 0.1223196930645258
$ws.Range("F12").Value = 1.712566634326478
$ws.Range("I12").Value = 1.063015125392624
$ws.Range("J12").Value = 0.1474484084356185
$ws.Range("K12").Value = 0.6229236479439351
$ws.Range("L12").Value = 0.3246081869731796
$ws.Range("N12").Value = 1.950497975728821
$ws.Range("O12").Value = 4.420407480748651
$ws.Range("B13").Value = 0.9086075759584276
$ws.Range("C13").Value = 0.07352882250192749
$ws.Range("D13").Value = 0.1369649222553591
$ws.Range("E13").Value = 0.122302941447586
$ws.Range("F13").Value = 1.712677008551054
$ws.Range("I13").Value = 1.063321534870667
$ws.Range("J13").Value = 0.1474840281780345
$ws.Range("K13").Value = 0.6205170074313457
$ws.Range("L13").Value = 0.3242759367403494
$ws.Range("N13").Value = 1.951476776734086
$ws.Range("O13").Value = 4.421166226057579
$ws.Range("B14").Value = 0.9001282654488421
$ws.Range("C14").Value = 0.0730591075211322
$ws.Range("D14").Value = 0.1364859221596788
$ws.Range("E14").Value = 0.1222492929232963
$ws.Range("F14").Value = 1.713055320173424
$ws.Range("I14").Value = 1.064331687948766
$ws.Range("J14").Value = 0.1476012796946904
$ws.Range("K14").Value = 0.612670447823092
$ws.Range("L14").Value = 0.3231952746906899
$ws.Range("N14").Value = 1.954683514660067
$ws.Range("O14").Value = 4.423689126727282
$ws.Range("B15").Value = 0.894938455183933
$ws.Range("C15").Value = 0.07277107184819442
$ws.Range("D15").Value = 0.1361933501617614
$ws.Range("E15").Value = 0.1222171868834963
$ws.Range("F15").Value = 1.713301156931536
$ws.Range("I15").Value = 1.064958801624414
$ws.Range("J15").Value = 0.1476739358265213
$ws.Range("K15").Value = 0.6078657664534717
$ws.Range("L15").Value = 0.3225355522062614
$ws.Range("N15").Value = 1.956658995630479
$ws.Range("O15").Value = 4.425271726103603
$ws.Range("B16").Value = 0.8652598680459676
$ws.Range("C16").Value = 0.07111554112935892
$ws.Range("D16").Value = 0.1345295045899491
$ws.Range("E16").Value = 0.1220447979176349
$ws.Range("F16").Value = 1.714926375846026
$ws.Range("I16").Value = 1.068680817716846
$ws.Range("J16").Value = 0.1481031280300478
$ws.Range("K16").Value = 0.5803565373995241
$ws.Range("L16").Value = 0.3187889940113706
$ws.Range("N16").Value = 1.968152291344683
$ws.Range("O16").Value = 4.434912178916534
$ws.Range("B17").Value = 0.8471091189431377
$ws.Range("C17").Value = 0.07009555834900993
$ws.Range("D17").Value = 0.1335202394502204
$ws.Range("E17").Value = 0.1219494158469416
$ws.Range("F17").Value = 1.716116691881282
$ws.Range("I17").Value = 1.071078683646157
$ws.Range("J17").Value = 0.1483778852202349
$ws.Range("K17").Value = 0.5635028854033237
$ws.Range("L17").Value = 0.3165210964288008
$ws.Range("N17").Value = 1.975357089230478
$ws.Range("O17").Value = 4.441336291064403
$ws.Range("B18").Value = 0.8366901730234702
$ws.Range("C18").Value = 0.06950729270349143
$ws.Range("D18").Value = 0.1329439680077655
$ws.Range("E18").Value = 0.1218983761833918
$ws.Range("F18").Value = 1.71687243521238
$ws.Range("I18").Value = 1.072499986548362
$ws.Range("J18").Value = 0.1485401359223459
$ws.Range("K18").Value = 0.5538175471360489
$ws.Range("L18").Value = 0.3152279142666714
$ws.Range("N18").Value = 1.979557708396705
$ws.Range("O18").Value = 4.445218815248467
$ws.Range("B19").Value = 0.8331661093080811
$ws.Range("C19").Value = 0.06930784216149277
$ws.Range("D19").Value = 0.1327495813107049
$ws.Range("E19").Value = 0.1218817522960762
$ws.Range("F19").Value = 1.717140535827497
$ws.Range("I19").Value = 1.072988450718391
$ws.Range("J19").Value = 0.1485957961277116
$ws.Range("K19").Value = 0.5505397274980339
$ws.Range("L19").Value = 0.3147920019941495
$ws.Range("N19").Value = 1.980989689031389
$ws.Range("O19").Value = 4.446565593524298
$ws.Range("B20").Value = 0.8490391416893033
$ws.Range("C20").Value = 0.07020430285135149
$ws.Range("D20").Value = 0.1336272400614718
$ws.Range("E20").Value = 0.1219591741273689
$ws.Range("F20").Value = 1.715982623160485
$ws.Range("I20").Value = 1.070819068654583
$ws.Range("J20").Value = 0.1483482004632393
$ws.Range("K20").Value = 0.5652961171446123
$ws.Range("L20").Value = 0.3167613544374177
$ws.Range("N20").Value = 1.974584267353115
$ws.Range("O20").Value = 4.44063302586099
$ws.Range("B21").Value = 0.9026180098221062
$ws.Range("C21").Value = 0.07319714127542909
$ws.Range("D21").Value = 0.1366264433836051
$ws.Range("E21").Value = 0.122264893459036
$ws.Range("F21").Value = 1.712941260759266
$ws.Range("I21").Value = 1.064033238225079
$ws.Range("J21").Value = 0.1475666658506967
$ws.Range("K21").Value = 0.6149748480657991
$ws.Range("L21").Value = 0.323512230143308
$ws.Range("N21").Value = 1.953739271147141
$ws.Range("O21").Value = 4.422940327466279
$ws.Range("B22").Value = 0.9378151957427292
$ws.Range("C22").Value = 0.07513875468755771
$ws.Range("D22").Value = 0.1386237802673804
$ws.Range("E22").Value = 0.1224985267141427
$ws.Range("F22").Value = 1.711585269107999
$ws.Range("I22").Value = 1.059972802920974
$ws.Range("J22").Value = 0.1470933285652904
$ws.Range("K22").Value = 0.6475133061062763
$ws.Range("L22").Value = 0.3280235446549398
$ws.Range("N22").Value = 1.940619902948224
$ws.Range("O22").Value = 4.413044662448954
$ws.Range("B23").Value = 0.9190135728763096
$ws.Range("C23").Value = 0.07410380181413245
$ws.Range("D23").Value = 0.1375543855439929
$ws.Range("E23").Value = 0.1223707467916952
$ws.Range("F23").Value = 1.712251256486326
$ws.Range("I23").Value = 1.062105692332587
$ws.Range("J23").Value = 0.147342537386205
$ws.Range("K23").Value = 0.630140687029126
$ws.Range("L23").Value = 0.3256067391708655
$ws.Range("N23").Value = 1.947575839347307
$ws.Range("O23").Value = 4.418173707699509
$ws.Range("B24").Value = 0.8481665280012578
$ws.Range("C24").Value = 0.07015514526634092
$ws.Range("D24").Value = 0.1335788527187702
$ws.Range("E24").Value = 0.1219547505802474
$ws.Range("F24").Value = 1.716043013161709
$ws.Range("I24").Value = 1.070936307461121
$ws.Range("J24").Value = 0.1483616075806999
$ws.Range("K24").Value = 0.5644853844057138
$ws.Range("L24").Value = 0.3166527005696196
$ws.Range("N24").Value = 1.974933477945435
$ws.Range("O24").Value = 4.440950382634725
$ws.Range("B25").Value = 0.7725480598081731
$ws.Range("C25").Value = 0.06583483573621862
$ws.Range("D25").Value = 0.1294526370647588
$ws.Range("E25").Value = 0.1216522541747
$ws.Range("F25").Value = 1.722852358993549
$ws.Range("I25").Value = 1.082071969781548
$ws.Range("J25").Value = 0.1496221594904874
$ws.Range("K25").Value = 0.4939905093785057
$ws.Range("L25").Value = 0.3074250612130385
$ws.Range("N25").Value = 2.006607566224706
$ws.Range("O25").Value = 4.472690388619185
